$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab
$ws.Name = "1주차"

# 2. Update header row 2: day-of-week labels (E2:K2), keep B2 "이름" as-is
$ws.Range("E2").Value = "월"
$ws.Range("F2").Value = "화"
$ws.Range("G2").Value = "수"
$ws.Range("H2").Value = "목"
$ws.Range("I2").Value = "금"
$ws.Range("J2").Value = "토"
$ws.Range("K2").Value = "일"

# 3. Update row 3: B3 becomes the author's own name, E3:K3 become day numbers 1-7
$ws.Range("B3").Value = "유도진"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 7

# 4. Select H7 to match final cursor position
$ws.Range("H7").Select()
